$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell Z1 is used as a staging area. Values are entered with a
# leading apostrophe so that numeric-looking text (e.g. "1.00", "36.70",
# thousand-dotted prices like "51.467.03") is kept as literal text instead
# of being auto-converted to a number by Excel. The value is then copied
# and pasted as "values only" into the real destination cell, so the
# destination keeps its original (default) cell style/format untouched.
$helper = $ws.Range("Z1")

function Set-TextValue($range, $text) {
    $helper.Value = "'" + $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") '51.467.03'
Set-TextValue $ws.Range("E2") '  +0.68%  '

Set-TextValue $ws.Range("D3") '2.980.33'
Set-TextValue $ws.Range("E3") '  +1.13%  '

Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.02%  '

Set-TextValue $ws.Range("D5") '381.52'
Set-TextValue $ws.Range("E5") '  +1.75%  '

Set-TextValue $ws.Range("D6") '103.52'
Set-TextValue $ws.Range("E6") '  +2.04%  '

Set-TextValue $ws.Range("D7") '0.545'
Set-TextValue $ws.Range("E7") '  +1.45%  '

Set-TextValue $ws.Range("D9") '0.592'
Set-TextValue $ws.Range("E9") '  +0.67%  '

Set-TextValue $ws.Range("D10") '36.70'
Set-TextValue $ws.Range("E10") '  +0.94%  '

Set-TextValue $ws.Range("E11") '  -0.80%  '

Set-TextValue $ws.Range("E12") '  +0.80%  '

Set-TextValue $ws.Range("D13") '3.444.88'
Set-TextValue $ws.Range("E13") '  +1.06%  '

Set-TextValue $ws.Range("D14") '18.43'
Set-TextValue $ws.Range("E14") '  +2.06%  '

Set-TextValue $ws.Range("D15") '7.78'
Set-TextValue $ws.Range("E15") '  +2.62%  '

Set-TextValue $ws.Range("D16") '2.982.50'
Set-TextValue $ws.Range("E16") '  -0.04%  '

Set-TextValue $ws.Range("D17") '11.18'
Set-TextValue $ws.Range("E17") '  +1.60%  '

Set-TextValue $ws.Range("D18") '0.995'
Set-TextValue $ws.Range("E18") '  -0.31%  '

Set-TextValue $ws.Range("D19") '51.461.12'
Set-TextValue $ws.Range("E19") '  +0.74%  '

Set-TextValue $ws.Range("D20") '3.08'
Set-TextValue $ws.Range("E20") '  -0.53%  '

Set-TextValue $ws.Range("D21") '12.62'
Set-TextValue $ws.Range("E21") '  +1.17%  '

Set-TextValue $ws.Range("D22") '0.0₃0962'
Set-TextValue $ws.Range("E22") '  +0.82%  '

Set-TextValue $ws.Range("D23") '70.44'
Set-TextValue $ws.Range("E23") '  +2.56%  '

Set-TextValue $ws.Range("D24") '267.42'
Set-TextValue $ws.Range("E24") '  +0.80%  '

Set-TextValue $ws.Range("D25") '3.22'
Set-TextValue $ws.Range("E25") '  +2.20%  '

Set-TextValue $ws.Range("D26") '7.86'
Set-TextValue $ws.Range("E26") '  -4.28%  '

Set-TextValue $ws.Range("D27") '7.32'
Set-TextValue $ws.Range("E27") '  -3.77%  '

Set-TextValue $ws.Range("D28") '0.169'
Set-TextValue $ws.Range("E28") '  +3.05%  '

Set-TextValue $ws.Range("D29") '0.999'
Set-TextValue $ws.Range("E29") '  -0.08%  '

Set-TextValue $ws.Range("D30") '26.08'
Set-TextValue $ws.Range("E30") '  +1.82%  '

Set-TextValue $ws.Range("D31") '0.109'
Set-TextValue $ws.Range("E31") '  -0.56%  '

Set-TextValue $ws.Range("D32") '10.33'
Set-TextValue $ws.Range("E32") '  +3.01%  '

Set-TextValue $ws.Range("D33") '34.62'
Set-TextValue $ws.Range("E33") '  +3.56%  '

Set-TextValue $ws.Range("D34") '51.57'
Set-TextValue $ws.Range("E34") '  +1.44%  '

Set-TextValue $ws.Range("E35") '  +0.60%  '

Set-TextValue $ws.Range("E36") '  -0.92%  '

Set-TextValue $ws.Range("E37") '  +0.10%  '

Set-TextValue $ws.Range("E38") '  +3.08%  '

Set-TextValue $ws.Range("D39") '16.79'
Set-TextValue $ws.Range("E39") '  +2.47%  '

Set-TextValue $ws.Range("E40") '  +1.46%  '

Set-TextValue $ws.Range("E41") '  +2.63%  '

Set-TextValue $ws.Range("E42") '  +2.06%  '

Set-TextValue $ws.Range("D43") '124.78'
Set-TextValue $ws.Range("E43") '  +3.33%  '

Set-TextValue $ws.Range("E44") '  +11.61%  '

Set-TextValue $ws.Range("D45") '21.41'
Set-TextValue $ws.Range("E45") '  +0.54%  '

Set-TextValue $ws.Range("E46") '  +0.06%  '

Set-TextValue $ws.Range("E47") '  +2.40%  '

Set-TextValue $ws.Range("E48") '  -1.12%  '

Set-TextValue $ws.Range("D49") '2.030.57'
Set-TextValue $ws.Range("E49") '  +1.89%  '

# Rows 50/51 swap places: row 50 (WOONetwork) <-> row 51 (BEAM), each with
# its own updated link/price/volume.
Set-TextValue $ws.Range("B50") 'BEAM'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
Set-TextValue $ws.Range("D50") '0.0334'
Set-TextValue $ws.Range("E50") '  +3.01%  '

Set-TextValue $ws.Range("B51") 'WOONetwork'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue $ws.Range("D51") '0.538'
Set-TextValue $ws.Range("E51") '  +16.37%  '

$helper.Clear()